$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: C710468 / NM1206B104K451CEGN / HUI JU capacitor
$ws.Range("A36").Value = "C710468"
$ws.Range("B36").Value = "NM1206B104K451CEGN"
$ws.Range("C36").Value = "HUI JU"
$ws.Range("D36").Value = 1206
$ws.Range("F36").Value = "100nF ±10% 450V X7R 1206 Multilayer Ceramic Capacitors MLCC - SMD/SMT RoHS"
$ws.Range("G36").Value = "yes"
$ws.Range("H36").Value = 20
$ws.Range("I36").Value = 20
$ws.Range("J36").Value = 0.0233
$ws.Range("K36").Value = 0.47
$ws.Range("L36").Value = "lcsc.com/product-detail/Multilayer-Ceramic-Capacitors-MLCC-SMD-SMT_HUI-JU-NM1206B104K451CEGN_C710468.html"

# Row 37: C357261 / TMPC1206HP-220MG-D / TAI-TECH inductor
$ws.Range("A37").Value = "C357261"
$ws.Range("B37").Value = "TMPC1206HP-220MG-D"
$ws.Range("C37").Value = "TAI-TECH"
$ws.Range("D37").Value = "SMD,13.5x12.5x5.7mm"
$ws.Range("F37").Value = "22uH ±20% 8A 34mΩ SMD,13.5x12.5x5.7mm Power Inductors RoHS"
$ws.Range("G37").Value = "yes"
$ws.Range("H37").Value = 2
$ws.Range("I37").Value = 1
$ws.Range("J37").Value = 0.6312
$ws.Range("K37").Value = 1.26
$ws.Range("L37").Value = "lcsc.com/product-detail/Power-Inductors_TAI-TECH-TMPC1206HP-220MG-D_C357261.html"

$ws.Rows("36:37").Select() | Out-Null
